$d = $word.ActiveDocument

# Locate the unique paragraph "...view_theme(position, vars): View HTML of position"
# and append " (deleted)" after the trailing "position", with "deleted" rendered
# italic + light-gray highlighted (split into "delete" + "d" runs, matching the
# target OOXML run boundaries).
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*View HTML of position*") {

        # " (" as a plain run, appended right before the paragraph mark.
        $p.Range.InsertAfter(" (")

        # "delete" as an italic + lightGray-highlighted run.
        $insertPos = $p.Range.End - 1
        $p.Range.InsertAfter("delete")
        $runDelete = $d.Range($insertPos, $insertPos + 6)
        $runDelete.Font.Italic = $true
        $runDelete.Font.HighlightColorIndex = 16

        # "d" as its own italic + lightGray-highlighted run (mirrors the diff's
        # two separate <w:r> elements for "delete" and "d").
        $insertPos2 = $p.Range.End - 1
        $p.Range.InsertAfter("d")
        $runD = $d.Range($insertPos2, $insertPos2 + 1)
        $runD.Font.Italic = $true
        $runD.Font.HighlightColorIndex = 16

        # Closing ")" as a plain run.
        $p.Range.InsertAfter(")")

        break
    }
}
